$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @("08-10-2021", "12-10-2021", "13-10-2021", "14-10-2021")
$data = @(
    @(14055, 20474, -6419),
    @(14046, 20354, -6308),
    @(14497, 20470, -5973),
    @(14815, 20263, -5448)
)

$startRow = 195
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    $cellA = $ws.Cells.Item($row, 1)
    $dateText = $dates[$i]

    # A few of these strings (e.g. "08-10-2021", "12-10-2021") look like
    # valid dates to Excel's auto-detection and would otherwise be silently
    # converted into a date serial number with date formatting applied.
    # Forcing a leading apostrophe makes Excel store the literal text
    # instead; ClearFormats() afterwards strips the "quote prefix" cell
    # style that the apostrophe trick adds, so the cell ends up as plain
    # text with no left-over formatting. Values that Excel would never
    # mistake for a date (e.g. "13-10-2021", where there is no 13th month)
    # are assigned directly, without needing the workaround.
    $testValue = $cellA.Parent.Application.IsNumber($dateText)
    $looksLikeDate = $false
    try {
        [void][datetime]::ParseExact($dateText, "dd-MM-yyyy", [System.Globalization.CultureInfo]::InvariantCulture)
        $looksLikeDate = $true
    } catch {
        $looksLikeDate = $false
    }

    if ($looksLikeDate) {
        $cellA.Value = "'" + $dateText
        $cellA.ClearFormats()
    } else {
        $cellA.Value = $dateText
    }

    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
}
